$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 8 (ano=2025) metrics with refreshed data
$ws.Range("C8").Value = 1241
$ws.Range("E8").Value = 1040
$ws.Range("G8").Value = 83.80338436744562
$ws.Range("H8").Value = 16.19661563255439
